# add on parse call back.
#
# Insert a new worksheet "user2" between the existing "user1" and "product"
# sheets. It gets a single "a,xxx" header/value pair (row2 header, row3
# data) like the other sheets, and becomes the active sheet/tab with A3
# selected. The selection remembered on "user1" is updated to A2:A3 and it
# is no longer the active tab.

$wb = $excel.ActiveWorkbook

# --- update the selection/active-cell remembered on the "user1" sheet ---
$user1 = $wb.Worksheets.Item("user1")
$user1.Activate()
$user1.Range("A2:A3").Select()

# --- insert the new "user2" sheet right after "user1" ---
# Duplicating "user1" (instead of Worksheets.Add) means the new sheet
# starts out with the same row height / formatting defaults used
# throughout the rest of this workbook.
$user1.Copy($null, $user1)
$user2 = $wb.Worksheets.Item(3)
$user2.Name = "user2"

# drop the copied data/columns, keeping just column A
$user2.Cells.Clear()
$user2.Columns("B:D").Delete()
$user2.Columns.Item(1).ColumnWidth = 43.15

# header row (row 2) mirrors the "name,type" convention used elsewhere
$user2.Range("A2").Value = "a,xxx"

# data row (row 3)
$user2.Range("A3").Value = 1
$user2.Range("A3").HorizontalAlignment = -4131

# make the new sheet the active tab with A3 selected
$user2.Activate()
$user2.Range("A3").Select()
